$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 is a new row that preserves the original (pre-edit) values of row 32.
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = "Femacal de La Calera"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44399
$ws.Range("D33").NumberFormat = $ws.Range("D32").NumberFormat
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 100112035
$ws.Range("G33").Value = "Bruselas (repollito)"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 38
$ws.Range("K33").Value = 22000
$ws.Range("L33").Value = 22000
$ws.Range("M33").Value = 22000
$ws.Range("N33").Value = '$/malla 15 kilos'
$ws.Range("O33").Value = "Provincia de Quillota"
$ws.Range("P33").Value = 1467
$ws.Range("Q33").Value = 15
$ws.Range("R33").Value = "Hortaliza"

# Row 32 is updated in-place with the new weekly report values.
$ws.Range("D32").Value = 44448
$ws.Range("J32").Value = 85
$ws.Range("K32").Value = 21000
$ws.Range("M32").Value = 21529
$ws.Range("P32").Value = 1435
